$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the formatting of existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the data rows: I column is always 1, J column mirrors column H
$data = @{
    2 = 2
    3 = 5
    4 = 6
    5 = 4
    6 = 4
    7 = 1
    8 = 3
    9 = 2
}

foreach ($row in $data.Keys) {
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $data[$row]
}
